$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Shanae", "Braun"),
    @("Johnson", "Ernser"),
    @("Jamey", "Weber"),
    @("Jame", "Senger"),
    @("Angelo", "Leffler"),
    @("TestName", "TestLastName")
)

$startRow = 7
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
